$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Zaccharie Risacher", "SF", "Atlanta Hawks"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder")
)

# Clear out old data rows (rows 2 through 19) first
$ws.Range("A2:C19").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
